$wb = $excel.ActiveWorkbook
$ws6 = $wb.Worksheets.Item("karamjeet")
$ws6.Activate()

# Header row: apply text number format (keeps existing bold + center)
$ws6.Range("A1:E1").NumberFormat = "@"

# Body area default text format (rows 2-8)
$ws6.Range("A2:E8").NumberFormat = "@"

# Row 2: re-enter id/dates as quoted text (values unchanged, now stored as text)
$ws6.Range("A2").Value = "'146"
$ws6.Range("D2").Value = "'2022-04-18"
$ws6.Range("E2").Value = "'2022-04-23"

# Row 3: new recruitment entry
$ws6.Range("A3").Value = "'17"
$ws6.Range("B3").Value = "recruitment_21"
$ws6.Range("C3").Value = "demo category"
$ws6.Range("D3").Value = "'2022-04-07"
$ws6.Range("E3").Value = "'2022-04-06"

$ws6.Range("E3").Select()
